{"js": "// Fix answer #5 for the Amdahl's law exercises: \"5. 0.3 seconds\" -> \"5. 28 seconds\"\nconst body = context.document.body;\n\n// The \"Answers:\" list item 5 is the last paragraph in the document body.\nlet paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\nlet target = paragraphs.items[paragraphs.items.length - 1];\n\n// Replace the wrong numeric answer (\"0.3\") with the corrected one (\"28\"),\n// scoped to that single paragraph so the identical \"0.3\" text in answer #4\n// is left untouched.\nlet hits = target.search(\"0.3\", { matchCase: true });\nhits.load(\"items\");\nawait context.sync();\nif (hits.items.length === 0) {\n  throw new Error('Could not find \"0.3\" in the last paragraph');\n}\nhits.items[0].insertText(\"28\", \"Replace\");\nawait context.sync();\n\n// The replace above leaves \"28\" merged into a single run together with the\n// rest of the paragraph's text. Force Word to split \"28\" back out into its\n// own run (so the run structure becomes \"5. \" / \"28\" / \" seconds\", mirroring\n// how answer #4 is already split into \"4. \" / \"0.3\" / \" seconds\") by\n// toggling a character property off and back onto its original value.\nparagraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\ntarget = paragraphs.items[paragraphs.items.length - 1];\nhits = target.search(\"28\", { matchCase: true });\nhits.load(\"items\");\nawait context.sync();\nhits.items[0].font.set({ bold: true });\nawait context.sync();\n\nparagraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\ntarget = paragraphs.items[paragraphs.items.length - 1];\nhits = target.search(\"28\", { matchCase: true });\nhits.load(\"items\");\nawait context.sync();\nhits.items[0].font.set({ bold: false });\nawait context.sync();\n", "ps1": "# Fix answer #5 for the Amdahl's law exercises: \"5. 0.3 seconds\" -> \"5. 28 seconds\"\n$d = $word.ActiveDocument\n\n# The \"Answers:\" list item 5 is the last paragraph in the document body.\n$para = $d.Paragraphs.Last\n\n# Replace the wrong numeric answer (\"0.3\") with the corrected one (\"28\"),\n# scoped to that single paragraph so the identical \"0.3\" in answer #4 is untouched.\n$range = $para.Range\n$find = $range.Find\n$find.Text = \"0.3\"\n$find.Execute() | Out-Null\nif ($find.Found) {\n    $range.Text = \"28\"\n}\n\n# The replace above leaves \"28\" merged into one run with the paragraph's\n# original text. Force Word to split \"28\" back out into its own run (so the\n# run structure matches \"5. \" / \"28\" / \" seconds\", mirroring how answer #4 is\n# already split into \"4. \" / \"0.3\" / \" seconds\") by nudging a character\n# property off and back onto its original value.\n$para2 = $d.Paragraphs.Last\n$range2 = $para2.Range\n$find2 = $range2.Find\n$find2.Text = \"28\"\n$find2.Execute() | Out-Null\nif ($find2.Found) {\n    $range2.Font.Bold = $true\n}\n\n$para3 = $d.Paragraphs.Last\n$range3 = $para3.Range\n$find3 = $range3.Find\n$find3.Text = \"28\"\n$find3.Execute() | Out-Null\nif ($find3.Found) {\n    $range3.Font.Bold = $false\n}\n"}
